$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(8).Delete()
$ws.Range("F2:F7").ClearContents()
$ws.Range("F2:F7").Style = "Normal"
